$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.056.66"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.406.95"
$ws.Range("E3").Value = "  -3.57%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "488.74"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.69%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "154.68"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.21%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.615"
$c.ClearFormats()
$ws.Range("E7").Value = "  +19.71%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "2.424.57"
$ws.Range("E9").Value = "  -3.42%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.26"
$c.ClearFormats()
$ws.Range("E10").Value = "  +9.39%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0998"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.53%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.333"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D14").Value = "2.826.92"
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("D15").Value = "57.089.34"
$ws.Range("E15").Value = "  +0.41%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "20.59"
$c.ClearFormats()
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "2.416.33"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("E19").Value = "  +4.09%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "325.12"
$c.ClearFormats()
$ws.Range("E20").Value = "  +1.18%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.96"
$c.ClearFormats()
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("E22").Value = "  -0.19%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.94"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.70%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "58.05"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.52%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.405"
$c.ClearFormats()
$ws.Range("E25").Value = "  -1.20%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "2.512.69"
$ws.Range("E28").Value = "  -3.71%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.25"
$c.ClearFormats()
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("D30").Value = "0.0₃0783"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("E31").Value = "  -0.01%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "150.49"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.85%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.53"
$c.ClearFormats()
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  -0.03%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.27"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.31%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.16"
$c.ClearFormats()
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.76"
$c.ClearFormats()
$ws.Range("E37").Value = "  -1.08%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.841"
$c.ClearFormats()
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("E39").Value = "  +8.67%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "34.18"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("E43").Value = "  -0.08%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.596"
$c.ClearFormats()
$ws.Range("E44").Value = "  -3.49%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "270.20"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -0.28%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.54"
$c.ClearFormats()
$ws.Range("E49").Value = "  -5.45%  "
$ws.Range("D50").Value = "1.880.14"
$ws.Range("E50").Value = "  -0.73%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "17.46"
$c.ClearFormats()
$ws.Range("E51").Value = "  -2.75%  "
